$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-03-12 Tuesday" "2024-03-13 Wednesday"

Replace-Text "954×8=" "651×9="
Replace-Text "394×9=" "839×2="
Replace-Text "987×6=" "232×3="
Replace-Text "486×2=" "823×8="
Replace-Text "968×2=" "852×2="
Replace-Text "726×4=" "521×2="
Replace-Text "315×5=" "285×5="
Replace-Text "599×4=" "957×5="
Replace-Text "488×8=" "934×2="
Replace-Text "737×4=" "502×6="
Replace-Text "346×3=" "244×2="
Replace-Text "498×7=" "368×3="
Replace-Text "164×7=" "221×3="
Replace-Text "164×4=" "461×4="
Replace-Text "230×6=" "443×9="
Replace-Text "329×3=" "678×4="
Replace-Text "334×2=" "618×7="
Replace-Text "771×4=" "258×7="
Replace-Text "479×7=" "173×9="
Replace-Text "517×6=" "157×7="
Replace-Text "630×4=" "465×9="
Replace-Text "787×4=" "967×8="
Replace-Text "648×7=" "690×8="
Replace-Text "286×9=" "947×5="
Replace-Text "354×4=" "182×8="
